$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "aman"
$ws.Range("B1").Value = "mohit"
$ws.Range("C1").Value = "tushar"
$ws.Range("D1").Value = "kashish"
$ws.Range("E1").Value = "yash"

$ws.Range("E1").Select()
